$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '94.918.79'
Set-TextValue $ws.Range('E2') '  -1.37%  '
Set-TextValue $ws.Range('D3') '3.554.91'
Set-TextValue $ws.Range('E3') '  -0.22%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '235.52'
Set-TextValue $ws.Range('E5') '  -1.85%  '
Set-TextValue $ws.Range('D6') '648.74'
Set-TextValue $ws.Range('E6') '  +1.78%  '
Set-TextValue $ws.Range('E7') '  -1.23%  '
Set-TextValue $ws.Range('D8') '0.397'
Set-TextValue $ws.Range('E8') '  -0.97%  '
Set-TextValue $ws.Range('E9') '  +0.04%  '
Set-TextValue $ws.Range('D10') '0.998'
Set-TextValue $ws.Range('E10') '  -2.29%  '
Set-TextValue $ws.Range('D11') '3.551.95'
Set-TextValue $ws.Range('E11') '  -0.25%  '
Set-TextValue $ws.Range('E12') '  +0.78%  '
Set-TextValue $ws.Range('D13') '42.20'
Set-TextValue $ws.Range('E13') '  -2.70%  '
Set-TextValue $ws.Range('D14') '6.46'
Set-TextValue $ws.Range('E14') '  +0.08%  '
Set-TextValue $ws.Range('D15') '4.260.81'
Set-TextValue $ws.Range('E15') '  +0.79%  '
Set-TextValue $ws.Range('D16') '94.868.57'
Set-TextValue $ws.Range('E16') '  -1.34%  '
Set-TextValue $ws.Range('E17') '  -0.79%  '
Set-TextValue $ws.Range('D18') '3.547.42'
Set-TextValue $ws.Range('E18') '  -0.46%  '
Set-TextValue $ws.Range('D19') '7.95'
Set-TextValue $ws.Range('E19') '  -0.73%  '
Set-TextValue $ws.Range('D20') '12.56'
Set-TextValue $ws.Range('E20') '  -4.18%  '
Set-TextValue $ws.Range('D21') '17.83'
Set-TextValue $ws.Range('E21') '  -1.59%  '
Set-TextValue $ws.Range('D22') '3.45'
Set-TextValue $ws.Range('E22') '  +0.55%  '
Set-TextValue $ws.Range('D23') '504.97'
Set-TextValue $ws.Range('E23') '  -2.22%  '
Set-TextValue $ws.Range('E24') '  -5.87%  '
Set-TextValue $ws.Range('E25') '  +1.10%  '
Set-TextValue $ws.Range('D26') '0.0000193'
Set-TextValue $ws.Range('E26') '  -0.82%  '
Set-TextValue $ws.Range('D27') '95.01'
Set-TextValue $ws.Range('E27') '  -1.99%  '
Set-TextValue $ws.Range('D28') '12.42'
Set-TextValue $ws.Range('E28') '  +0.54%  '
Set-TextValue $ws.Range('D29') '3.745.30'
Set-TextValue $ws.Range('E29') '  -0.16%  '
Set-TextValue $ws.Range('D30') '3.00'
Set-TextValue $ws.Range('E30') '  -3.56%  '
Set-TextValue $ws.Range('D31') '11.37'
Set-TextValue $ws.Range('E31') '  -1.63%  '
Set-TextValue $ws.Range('E32') '  -3.50%  '
Set-TextValue $ws.Range('E33') '  +0.01%  '
Set-TextValue $ws.Range('D34') '0.999'
Set-TextValue $ws.Range('E34') '  -0.88%  '
Set-TextValue $ws.Range('D35') '0.176'
Set-TextValue $ws.Range('E35') '  -3.72%  '
Set-TextValue $ws.Range('D36') '31.62'
Set-TextValue $ws.Range('E36') '  +4.69%  '
Set-TextValue $ws.Range('D37') '0.556'
Set-TextValue $ws.Range('E37') '  -1.58%  '
Set-TextValue $ws.Range('D38') '8.46'
Set-TextValue $ws.Range('E38') '  +7.02%  '
Set-TextValue $ws.Range('D39') '1.61'
Set-TextValue $ws.Range('E39') '  +6.64%  '
Set-TextValue $ws.Range('D40') '583.73'
Set-TextValue $ws.Range('E40') '  +0.56%  '
Set-TextValue $ws.Range('E41') '  +0.03%  '
Set-TextValue $ws.Range('E42') '  -1.35%  '
Set-TextValue $ws.Range('D43') '0.899'
Set-TextValue $ws.Range('E43') '  -2.45%  '
Set-TextValue $ws.Range('D44') '1.75'
Set-TextValue $ws.Range('E44') '  -0.27%  '
Set-TextValue $ws.Range('D45') '2.28'
Set-TextValue $ws.Range('E45') '  +4.48%  '
Set-TextValue $ws.Range('D46') '5.65'
Set-TextValue $ws.Range('E46') '  +0.59%  '
Set-TextValue $ws.Range('D47') '23.37'
Set-TextValue $ws.Range('E47') '  -2.11%  '
Set-TextValue $ws.Range('B48') 'VeChain'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D48') '0.0411'
Set-TextValue $ws.Range('E48') '  -5.23%  '
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '33.32'
Set-TextValue $ws.Range('E49') '  +32.10%  '
Set-TextValue $ws.Range('D50') '3.57'
Set-TextValue $ws.Range('E50') '  +0.69%  '
Set-TextValue $ws.Range('D51') '53.21'
Set-TextValue $ws.Range('E51') '  -1.36%  '
